$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 04:05"

# Pakistan (row 22) - updated counts
$ws.Range("B22").Value = 34336
$ws.Range("C22").Value = 1662
$ws.Range("D22").Value = 8812
$ws.Range("E22").Value = 24787
$ws.Range("F22").Value = 111
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 737

# Australia (row 53) - updated counts
$ws.Range("B53").Value = 6980
$ws.Range("C53").Value = 16
$ws.Range("D53").Value = 6270
$ws.Range("E53").Value = 612
$ws.Range("F53").Value = 18
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 98

# Bolivia overtakes Irak in the ranking: row 70 becomes Bolivia (updated data),
# row 71 becomes Irak (keeps Irak's prior data unchanged)
$ws.Range("A70").Value = "Bolivia"
$ws.Range("B70").Value = 2964
$ws.Range("C70").Value = 133
$ws.Range("D70").Value = 313
$ws.Range("E70").Value = 2523
$ws.Range("F70").Value = 3
$ws.Range("G70").Value = 6
$ws.Range("H70").Value = 128

$ws.Range("A71").Value = "Irak"
$ws.Range("B71").Value = 2913
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 1903
$ws.Range("E71").Value = 898
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 112

# Guatemala overtakes Somalia in the ranking: row 93 becomes Guatemala (updated data),
# row 94 becomes Somalia (keeps Somalia's prior data unchanged)
$ws.Range("A93").Value = "Guatemala"
$ws.Range("B93").Value = 1199
$ws.Range("C93").Value = 85
$ws.Range("D93").Value = 120
$ws.Range("E93").Value = 1052
$ws.Range("F93").Value = 5
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 27

$ws.Range("A94").Value = "Somalia"
$ws.Range("B94").Value = 1170
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 126
$ws.Range("E94").Value = 992
$ws.Range("F94").Value = 2
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 52
